# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E22) is re-sorted into chronological order
# (oldest -> newest): 2212, 2301, 2302, 2303, 2304, 2305, 2306
# (it used to run newest -> oldest: 2306, 2305, 2304, 2303, 2302, 2301, 2212).
#
# The "Valor Mora" column (F16:F22) values travel together with their
# period label, so only the two periods that actually swapped rows
# (2212 and 2306) end up with different F values than before:
#   - period 2212 moves from row 22 to row 16, taking its value (40000) along
#   - period 2306 moves from row 16 to row 22, taking its value (29333) along

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New chronological order of periods for E16:E22
$periods = @("2212", "2301", "2302", "2303", "2304", "2305", "2306")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# Valor Mora follows the period label: 2212 -> 40000, 2306 -> 29333
$ws.Range("F16").Value = 40000
$ws.Range("F22").Value = 29333
